# Update the multiplication problems in the single table, addressed by
# (row, column) so that the textually-identical "750x6=" value appearing
# both as a source and a destination string cannot cause a double match.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1; Col = 1; Old = "374×5="; New = "626×9=" }
    @{ Row = 1; Col = 2; Old = "996×9="; New = "485×6=" }
    @{ Row = 1; Col = 3; Old = "277×4="; New = "304×5=" }
    @{ Row = 1; Col = 4; Old = "177×3="; New = "656×8=" }
    @{ Row = 1; Col = 5; Old = "409×2="; New = "266×3=" }
    @{ Row = 5; Col = 1; Old = "117×2="; New = "782×5=" }
    @{ Row = 5; Col = 2; Old = "870×7="; New = "177×6=" }
    @{ Row = 5; Col = 3; Old = "442×2="; New = "376×7=" }
    @{ Row = 5; Col = 4; Old = "632×5="; New = "715×2=" }
    @{ Row = 5; Col = 5; Old = "620×4="; New = "613×7=" }
    @{ Row = 10; Col = 1; Old = "666×2="; New = "750×6=" }
    @{ Row = 10; Col = 2; Old = "809×7="; New = "542×3=" }
    @{ Row = 10; Col = 3; Old = "475×4="; New = "126×2=" }
    @{ Row = 10; Col = 4; Old = "454×3="; New = "537×5=" }
    @{ Row = 10; Col = 5; Old = "985×9="; New = "559×6=" }
    @{ Row = 15; Col = 1; Old = "174×9="; New = "832×6=" }
    @{ Row = 15; Col = 2; Old = "205×9="; New = "500×4=" }
    @{ Row = 15; Col = 3; Old = "197×3="; New = "147×7=" }
    @{ Row = 15; Col = 4; Old = "477×2="; New = "952×8=" }
    @{ Row = 15; Col = 5; Old = "443×9="; New = "457×3=" }
    @{ Row = 20; Col = 1; Old = "390×5="; New = "904×9=" }
    @{ Row = 20; Col = 2; Old = "750×6="; New = "114×4=" }
    @{ Row = 20; Col = 3; Old = "660×3="; New = "747×3=" }
    @{ Row = 20; Col = 4; Old = "718×2="; New = "994×5=" }
    @{ Row = 20; Col = 5; Old = "843×2="; New = "172×4=" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $range = $cell.Range
    # Cell.Range.Text includes the trailing cell-mark characters; trim them
    # off when comparing against the expected plain-text value.
    $current = $range.Text.TrimEnd([char]7, [char]13, [char]10)
    if ($current -ne $u.Old) {
        throw "Unexpected cell text at row $($u.Row) col $($u.Col): expected $($u.Old) but found $current"
    }
    $range.Text = $u.New
}

Write-Output "Updated $($updates.Count) cells"
